$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:U2").Value = 83.59375
